# Update the "lower bits of COPPER ID" box on slide 7:
#   (11-0)    lower bits of COPPER ID  :  12bit (1024)
# becomes
#   (9-0)    lower bits of COPPER ID  :  10bit (1024)
#
# The node-id bit range "11-0" / "12bit" is re-typed so the changed
# digits land in their own runs, matching how the original author
# retyped the text in PowerPoint (one run boundary per edited digit).
#
# Note: this runtime's TextRange.LanguageID setter does not persist to
# the saved OOXML (and, worse, corrupts unrelated runs in the shape), so
# it is deliberately not used here - only Characters(start,len).Text
# assignments are used, which reliably split runs at the given
# boundaries.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Locate the shape (a rectangle callout) that holds the text we need to
# edit - its shape Id is 5 ("正方形/長方形 4"). Look it up by Id rather than
# a bare positional index so the script stays correct even if shapes are
# reordered.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 5) {
        $targetShape = $candidate
        break
    }
}

$tr = $targetShape.TextFrame.TextRange

# ---------------------------------------------------------------------
# The shape's second paragraph reads (before any edits):
#   "(11-0)    lower bits of COPPER ID  :  12bit (1024)"
# Character columns (1-based) within the whole shape's TextRange:
#   49        "("
#   50-51     "11"
#   52        "-"
#   53        "0"
#   54        ")"
#   55-58     "    "
#   ...
#   84-86     ":  "
#   87        "1"
#   88        "2"
#   89-92     "bit "
#   93        "("
#   94-97     "1024"
#   98        ")"
#
# Edit the higher-numbered region first ("12bit" -> "10bit") so the
# still-to-come shrink in the lower region ("11-0" -> "9-0", 9 chars ->
# 8 chars) can't shift the offsets used below.
# ---------------------------------------------------------------------

# ":  12bit (" -> ":  " / "1" / "0" / "bit ("
$tr.Characters(87, 1).Text = "1"
$tr.Characters(88, 1).Text = "0"
$tr.Characters(89, 5).Text = "bit ("

# "11-0)    " -> "9" / "-0" / ")    "
$tr.Characters(50, 2).Text = "9"
$tr.Characters(51, 2).Text = "-0"
$tr.Characters(53, 5).Text = ")    "
